$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Forwarder" header column (D) next to the existing
# AwbNo / OrderNo / BagNo headers, matching the style (bold font)
# used by the other headers but highlighted with a yellow fill so
# it stands out on the printed label.
$ws.Range("D1").Value = "Forwarder"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Interior.Color = 65535

# Size the new column to fit its header text.
$ws.Columns.Item(4).AutoFit()

# Leave the selection where the user last clicked while editing.
[void]$ws.Range("E4").Select()
